$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "Passed" -> "Status" ---
$ws.Range("C1").Value = "Status"

# --- Row 2: function renamed from decimalToBinary -> binaryToDecimal ---
$ws.Range("A2").Value = "binaryToDecimal"

# --- "Yes" -> "Passed" for every data row in column C (rows 2-8) ---
$ws.Range("C2").Value = "Passed"
$ws.Range("C3").Value = "Passed"
$ws.Range("C4").Value = "Passed"
$ws.Range("C5").Value = "Passed"
$ws.Range("C6").Value = "Passed"
$ws.Range("C7").Value = "Passed"
$ws.Range("C8").Value = "Passed"

# --- Append a new row 9 (userAnswer), copying formatting from row 8 ---
$ws.Range("A8:E8").Copy($ws.Range("A9:E9"))

$ws.Range("A9").Value = "userAnswer"
$ws.Range("B9").Value = "Prompts the player for an answer"
$ws.Range("C9").Value = "Passed"
$ws.Range("D9").Value = "If the answer is wrong, asks the user for a new answer, otherwise stops the program"
$ws.Range("E9").Value = "If the answer is wrong, asks the user for a new answer, otherwise stops the program"

# --- Row height adjustments ---
$ws.Rows.Item(8).RowHeight = 18.75
$ws.Rows.Item(9).RowHeight = 19.5

# --- Column width adjustments: split the old D:E 43-wide pair into distinct widths ---
$ws.Columns.Item(4).ColumnWidth = 74.16666666666667
$ws.Columns.Item(5).ColumnWidth = 74

# --- Update the selected / active cell ---
$ws.Range("E14").Select()
